$wb = $excel.ActiveWorkbook

# Map: row -> new "want to go" count (column F), per worksheet, as updated
# by the gh-pages data refresh (commit 456a3b4).

# sheet1: sheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 65
$ws.Cells.Item(5, 6).Value = 7564
$ws.Cells.Item(7, 6).Value = 7761
$ws.Cells.Item(9, 6).Value = 42
$ws.Cells.Item(11, 6).Value = 6427
$ws.Cells.Item(12, 6).Value = 3317
$ws.Cells.Item(14, 6).Value = 3681
$ws.Cells.Item(15, 6).Value = 35
$ws.Cells.Item(16, 6).Value = 31
$ws.Cells.Item(17, 6).Value = 32
$ws.Cells.Item(18, 6).Value = 48
$ws.Cells.Item(19, 6).Value = 23
$ws.Cells.Item(23, 6).Value = 303
$ws.Cells.Item(24, 6).Value = 3740
$ws.Cells.Item(26, 6).Value = 355
$ws.Cells.Item(27, 6).Value = 948
$ws.Cells.Item(28, 6).Value = 273
$ws.Cells.Item(29, 6).Value = 1399
$ws.Cells.Item(31, 6).Value = 40
$ws.Cells.Item(32, 6).Value = 2698
$ws.Cells.Item(33, 6).Value = 1688
$ws.Cells.Item(34, 6).Value = 26
$ws.Cells.Item(37, 6).Value = 3501
$ws.Cells.Item(38, 6).Value = 254
$ws.Cells.Item(41, 6).Value = 912
$ws.Cells.Item(42, 6).Value = 513
$ws.Cells.Item(43, 6).Value = 1357
$ws.Cells.Item(44, 6).Value = 239
$ws.Cells.Item(46, 6).Value = 617

# sheet2: sheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(3, 6).Value = 238
$ws.Cells.Item(5, 6).Value = 22
$ws.Cells.Item(6, 6).Value = 42
$ws.Cells.Item(8, 6).Value = 37
$ws.Cells.Item(9, 6).Value = 54
$ws.Cells.Item(14, 6).Value = 82

# sheet3: sheet index 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 128

# sheet4: sheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 6).Value = 65
$ws.Cells.Item(5, 6).Value = 128
$ws.Cells.Item(7, 6).Value = 238
$ws.Cells.Item(9, 6).Value = 42
$ws.Cells.Item(10, 6).Value = 7564
$ws.Cells.Item(11, 6).Value = 7761
$ws.Cells.Item(14, 6).Value = 6427
$ws.Cells.Item(15, 6).Value = 3317
$ws.Cells.Item(16, 6).Value = 3681
$ws.Cells.Item(17, 6).Value = 31
$ws.Cells.Item(18, 6).Value = 48
$ws.Cells.Item(19, 6).Value = 23
$ws.Cells.Item(21, 6).Value = 37
$ws.Cells.Item(24, 6).Value = 3740
$ws.Cells.Item(28, 6).Value = 355
$ws.Cells.Item(29, 6).Value = 948
$ws.Cells.Item(30, 6).Value = 273
$ws.Cells.Item(31, 6).Value = 1399
$ws.Cells.Item(33, 6).Value = 40
$ws.Cells.Item(34, 6).Value = 2698
$ws.Cells.Item(35, 6).Value = 1688
$ws.Cells.Item(36, 6).Value = 26
$ws.Cells.Item(39, 6).Value = 82
$ws.Cells.Item(40, 6).Value = 3501
$ws.Cells.Item(41, 6).Value = 255
$ws.Cells.Item(44, 6).Value = 912
$ws.Cells.Item(45, 6).Value = 513
$ws.Cells.Item(46, 6).Value = 1357
$ws.Cells.Item(49, 6).Value = 617
